$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1: "audioFalse" -> "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Update C2 and C3 to the new shared value "train1P2"
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
